$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new task entry for row 19 (sr no 18, date 45436)
$ws.Range("C19").Value = "creating a profile page ui"
$ws.Range("E19").Value = "1 day"

# Update the active view/selection to reflect where the user left off
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("E19").Select()
